$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 2-13 with new TPM-derived values ---
$ws.Range("G2").Value2 = 10.100659
$ws.Range("H2").Value2 = 30.301977
$ws.Range("I2").Value2 = 0.3328245842863797
$ws.Range("J2").Value2 = 0.3328245842863797
$ws.Range("M2").Value2 = 2.885873333333334
$ws.Range("N2").Value2 = 8.657620000000001
$ws.Range("O2").Value2 = 0.3070415651026022
$ws.Range("P2").Value2 = 0.3070415651026022
$ws.Range("Q2").Value2 = 29.14922245719334
$ws.Range("R2").Value2 = 262.34300211474
$ws.Range("S2").Value2 = 0.102190981263913
$ws.Range("T2").Value2 = 0.102190981263913
$ws.Range("G3").Value2 = 10.100659
$ws.Range("H3").Value2 = 30.301977
$ws.Range("I3").Value2 = 0.3328245842863797
$ws.Range("J3").Value2 = 0.3328245842863797
$ws.Range("O3").Value2 = 0.3368406220840099
$ws.Range("P3").Value2 = 0.3368406220840099
$ws.Range("Q3").Value2 = 31.97821839679967
$ws.Range("R3").Value2 = 287.803965571197
$ws.Range("S3").Value2 = 0.1121088400158761
$ws.Range("T3").Value2 = 0.1121088400158761
$ws.Range("G4").Value2 = 10.100659
$ws.Range("H4").Value2 = 30.301977
$ws.Range("I4").Value2 = 0.3328245842863797
$ws.Range("J4").Value2 = 0.3328245842863797
$ws.Range("M4").Value2 = 3.327024333333334
$ws.Range("N4").Value2 = 9.981073
$ws.Range("O4").Value2 = 0.3539776838580724
$ws.Range("P4").Value2 = 0.3539776838580724
$ws.Range("Q4").Value2 = 33.60513827570234
$ws.Range("R4").Value2 = 302.446244481321
$ws.Range("S4").Value2 = 0.1178124754767185
$ws.Range("T4").Value2 = 0.1178124754767185
$ws.Range("A5").Value2 = "ECs"
$ws.Range("D5").Value2 = "Resolving-Mac"
$ws.Range("G5").Value2 = 10.100659
$ws.Range("H5").Value2 = 30.301977
$ws.Range("I5").Value2 = 0.3328245842863797
$ws.Range("J5").Value2 = 0.3328245842863797
$ws.Range("K5").Value2 = 1
$ws.Range("L5").Value2 = 0.3333333333333333
$ws.Range("M5").Value2 = 0.020115
$ws.Range("N5").Value2 = 0.060345
$ws.Range("O5").Value2 = 0.002140128955315263
$ws.Range("P5").Value2 = 0.002140128955315263
$ws.Range("Q5").Value2 = 0.203174755785
$ws.Range("R5").Value2 = 1.828572802065
$ws.Range("S5").Value2 = 0.0007122875298720464
$ws.Range("T5").Value2 = 0.0007122875298720466
$ws.Range("D6").Value2 = "ECs"
$ws.Range("I6").Value2 = 0.4180918757349671
$ws.Range("J6").Value2 = 0.4180918757349671
$ws.Range("M6").Value2 = 2.885873333333334
$ws.Range("N6").Value2 = 8.657620000000001
$ws.Range("O6").Value2 = 0.3070415651026022
$ws.Range("P6").Value2 = 0.3070415651026022
$ws.Range("Q6").Value2 = 36.61704594170667
$ws.Range("R6").Value2 = 329.5534134753601
$ws.Range("S6").Value2 = 0.128371583882347
$ws.Range("T6").Value2 = 0.128371583882347
$ws.Range("D7").Value2 = "FAPs"
$ws.Range("I7").Value2 = 0.4180918757349671
$ws.Range("J7").Value2 = 0.4180918757349671
$ws.Range("M7").Value2 = 3.165953666666667
$ws.Range("N7").Value2 = 9.497861
$ws.Range("O7").Value2 = 0.3368406220840099
$ws.Range("P7").Value2 = 0.3368406220840099
$ws.Range("Q7").Value2 = 40.17081052124534
$ws.Range("R7").Value2 = 361.537294691208
$ws.Range("S7").Value2 = 0.1408303275108369
$ws.Range("T7").Value2 = 0.1408303275108369
$ws.Range("A8").Value2 = "FAPs"
$ws.Range("D8").Value2 = "MuSCs"
$ws.Range("G8").Value2 = 12.688376
$ws.Range("H8").Value2 = 38.065128
$ws.Range("I8").Value2 = 0.4180918757349671
$ws.Range("J8").Value2 = 0.4180918757349671
$ws.Range("M8").Value2 = 3.327024333333334
$ws.Range("N8").Value2 = 9.981073
$ws.Range("O8").Value2 = 0.3539776838580724
$ws.Range("P8").Value2 = 0.3539776838580724
$ws.Range("Q8").Value2 = 42.21453570248267
$ws.Range("R8").Value2 = 379.930821322344
$ws.Range("S8").Value2 = 0.1479951938125407
$ws.Range("T8").Value2 = 0.1479951938125407
$ws.Range("A9").Value2 = "FAPs"
$ws.Range("D9").Value2 = "Resolving-Mac"
$ws.Range("G9").Value2 = 12.688376
$ws.Range("H9").Value2 = 38.065128
$ws.Range("I9").Value2 = 0.4180918757349671
$ws.Range("J9").Value2 = 0.4180918757349671
$ws.Range("K9").Value2 = 1
$ws.Range("L9").Value2 = 0.3333333333333333
$ws.Range("M9").Value2 = 0.020115
$ws.Range("N9").Value2 = 0.060345
$ws.Range("O9").Value2 = 0.002140128955315263
$ws.Range("P9").Value2 = 0.002140128955315263
$ws.Range("Q9").Value2 = 0.25522668324
$ws.Range("R9").Value2 = 2.29704014916
$ws.Range("S9").Value2 = 0.0008947705292424739
$ws.Range("T9").Value2 = 0.000894770529242474
$ws.Range("D10").Value2 = "ECs"
$ws.Range("G10").Value2 = 4.721016333333334
$ws.Range("H10").Value2 = 14.163049
$ws.Range("I10").Value2 = 0.1555611667071302
$ws.Range("J10").Value2 = 0.1555611667071302
$ws.Range("M10").Value2 = 2.885873333333334
$ws.Range("N10").Value2 = 8.657620000000001
$ws.Range("O10").Value2 = 0.3070415651026022
$ws.Range("P10").Value2 = 0.3070415651026022
$ws.Range("Q10").Value2 = 13.62425514259778
$ws.Range("R10").Value2 = 122.61829628338
$ws.Range("S10").Value2 = 0.04776374409494409
$ws.Range("T10").Value2 = 0.04776374409494408
$ws.Range("A11").Value2 = "MuSCs"
$ws.Range("D11").Value2 = "FAPs"
$ws.Range("G11").Value2 = 4.721016333333334
$ws.Range("H11").Value2 = 14.163049
$ws.Range("I11").Value2 = 0.1555611667071302
$ws.Range("J11").Value2 = 0.1555611667071302
$ws.Range("M11").Value2 = 3.165953666666667
$ws.Range("N11").Value2 = 9.497861
$ws.Range("O11").Value2 = 0.3368406220840099
$ws.Range("P11").Value2 = 0.3368406220840099
$ws.Range("Q11").Value2 = 14.94651897090989
$ws.Range("R11").Value2 = 134.518670738189
$ws.Range("S11").Value2 = 0.05239932016574413
$ws.Range("T11").Value2 = 0.05239932016574412
$ws.Range("A12").Value2 = "MuSCs"
$ws.Range("D12").Value2 = "MuSCs"
$ws.Range("G12").Value2 = 4.721016333333334
$ws.Range("H12").Value2 = 14.163049
$ws.Range("I12").Value2 = 0.1555611667071302
$ws.Range("J12").Value2 = 0.1555611667071302
$ws.Range("M12").Value2 = 3.327024333333334
$ws.Range("N12").Value2 = 9.981073
$ws.Range("O12").Value2 = 0.3539776838580724
$ws.Range("P12").Value2 = 0.3539776838580724
$ws.Range("Q12").Value2 = 15.70693621906412
$ws.Range("R12").Value2 = 141.362425971577
$ws.Range("S12").Value2 = 0.05506518148924945
$ws.Range("T12").Value2 = 0.05506518148924944
$ws.Range("A13").Value2 = "MuSCs"
$ws.Range("D13").Value2 = "Resolving-Mac"
$ws.Range("G13").Value2 = 4.721016333333334
$ws.Range("H13").Value2 = 14.163049
$ws.Range("I13").Value2 = 0.1555611667071302
$ws.Range("J13").Value2 = 0.1555611667071302
$ws.Range("K13").Value2 = 1
$ws.Range("L13").Value2 = 0.3333333333333333
$ws.Range("M13").Value2 = 0.020115
$ws.Range("N13").Value2 = 0.060345
$ws.Range("O13").Value2 = 0.002140128955315263
$ws.Range("P13").Value2 = 0.002140128955315263
$ws.Range("Q13").Value2 = 0.09496324354500002
$ws.Range("R13").Value2 = 0.8546691919050001
$ws.Range("S13").Value2 = 0.0003329209571925541
$ws.Range("T13").Value2 = 0.0003329209571925541

# --- Add new rows 14-17 ---
$ws.Range("A14").Value2 = "Resolving-Mac"
$ws.Range("B14").Value2 = "Hbegf"
$ws.Range("C14").Value2 = "Erbb2"
$ws.Range("D14").Value2 = "ECs"
$ws.Range("E14").Value2 = 3
$ws.Range("F14").Value2 = 1
$ws.Range("G14").Value2 = 2.838244666666667
$ws.Range("H14").Value2 = 8.514734000000001
$ws.Range("I14").Value2 = 0.09352237327152295
$ws.Range("J14").Value2 = 0.09352237327152294
$ws.Range("K14").Value2 = 3
$ws.Range("L14").Value2 = 1
$ws.Range("M14").Value2 = 2.885873333333334
$ws.Range("N14").Value2 = 8.657620000000001
$ws.Range("O14").Value2 = 0.3070415651026022
$ws.Range("P14").Value2 = 0.3070415651026022
$ws.Range("Q14").Value2 = 8.190814597008892
$ws.Range("R14").Value2 = 73.71733137308001
$ws.Range("S14").Value2 = 0.02871525586139818
$ws.Range("T14").Value2 = 0.02871525586139817

$ws.Range("A15").Value2 = "Resolving-Mac"
$ws.Range("B15").Value2 = "Hbegf"
$ws.Range("C15").Value2 = "Erbb2"
$ws.Range("D15").Value2 = "FAPs"
$ws.Range("E15").Value2 = 3
$ws.Range("F15").Value2 = 1
$ws.Range("G15").Value2 = 2.838244666666667
$ws.Range("H15").Value2 = 8.514734000000001
$ws.Range("I15").Value2 = 0.09352237327152295
$ws.Range("J15").Value2 = 0.09352237327152294
$ws.Range("K15").Value2 = 3
$ws.Range("L15").Value2 = 1
$ws.Range("M15").Value2 = 3.165953666666667
$ws.Range("N15").Value2 = 9.497861
$ws.Range("O15").Value2 = 0.3368406220840099
$ws.Range("P15").Value2 = 0.3368406220840099
$ws.Range("Q15").Value2 = 8.985751109330446
$ws.Range("R15").Value2 = 80.87175998397402
$ws.Range("S15").Value2 = 0.03150213439155278
$ws.Range("T15").Value2 = 0.03150213439155277

$ws.Range("A16").Value2 = "Resolving-Mac"
$ws.Range("B16").Value2 = "Hbegf"
$ws.Range("C16").Value2 = "Erbb2"
$ws.Range("D16").Value2 = "MuSCs"
$ws.Range("E16").Value2 = 3
$ws.Range("F16").Value2 = 1
$ws.Range("G16").Value2 = 2.838244666666667
$ws.Range("H16").Value2 = 8.514734000000001
$ws.Range("I16").Value2 = 0.09352237327152295
$ws.Range("J16").Value2 = 0.09352237327152294
$ws.Range("K16").Value2 = 3
$ws.Range("L16").Value2 = 1
$ws.Range("M16").Value2 = 3.327024333333334
$ws.Range("N16").Value2 = 9.981073
$ws.Range("O16").Value2 = 0.3539776838580724
$ws.Range("P16").Value2 = 0.3539776838580724
$ws.Range("Q16").Value2 = 9.442909069953556
$ws.Range("R16").Value2 = 84.98618162958201
$ws.Range("S16").Value2 = 0.03310483307956379
$ws.Range("T16").Value2 = 0.03310483307956379

$ws.Range("A17").Value2 = "Resolving-Mac"
$ws.Range("B17").Value2 = "Hbegf"
$ws.Range("C17").Value2 = "Erbb2"
$ws.Range("D17").Value2 = "Resolving-Mac"
$ws.Range("E17").Value2 = 3
$ws.Range("F17").Value2 = 1
$ws.Range("G17").Value2 = 2.838244666666667
$ws.Range("H17").Value2 = 8.514734000000001
$ws.Range("I17").Value2 = 0.09352237327152295
$ws.Range("J17").Value2 = 0.09352237327152294
$ws.Range("K17").Value2 = 1
$ws.Range("L17").Value2 = 0.3333333333333333
$ws.Range("M17").Value2 = 0.020115
$ws.Range("N17").Value2 = 0.060345
$ws.Range("O17").Value2 = 0.002140128955315263
$ws.Range("P17").Value2 = 0.002140128955315263
$ws.Range("Q17").Value2 = 0.05709129147000001
$ws.Range("R17").Value2 = 0.51382162323
$ws.Range("S17").Value2 = 0.0002001499390081885
$ws.Range("T17").Value2 = 0.0002001499390081885

